$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequence")

# Insert two new blank rows before row 5 (pushes existing rows 5+ down to 7+)
$ws.Rows("5:6").Insert()

# Give the new banner row (B5:L5) the same look as the existing explanation
# banner (B10:L10, style index 18: wrap text + left aligned) without creating
# a brand new style table entry.
$ws.Range("B10:L10").Copy()
$ws.Range("B5:L5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New explanatory paragraph text for the merged banner cell B5:L5
$newText = @"
This sheet is used to document an inferred gene sequence published by IARC, based on information contained
within one or more submissions to the Committee. It may be re-issued if new information comes to light (for
example, if supporting information is contained in further submissions to IARC). In such cases, the description_id 
will never change, and can therefore be used to associate revisions of the same inferred sequence.
"@

$ws.Range("B5").Value = $newText
$ws.Range("B5:L5").Merge()
$ws.Rows(5).RowHeight = 56.25

# Row 6 stays blank (formatting already copied down by the row Insert above)

# Row 7 now holds "Please complete details on all tabs." - give it the explicit row height
$ws.Rows(7).RowHeight = 16.5

# Keep selection / view consistent with the reviewed workbook
$ws.Application.GoTo($ws.Range("A3"))
$ws.Range("P5").Select()
